# Normalise the "Recorded By" column (G): entries that were recorded by a
# mix of a named account and the automated "System" account are listed as
# comma-separated values. Re-order each such list so that "System"/"system"
# is promoted to the front, e.g. "user@example.com, System" becomes
# "System, user@example.com".
#
# A cell is only touched when it contains more than one comma-separated
# value AND the last value in that list is "System" (case-insensitive);
# everything else (single-value cells, or lists that already start with
# "System" and have no trailing "System" entry) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Text

    if ($val -eq $null -or $val -eq "") { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $lastPart = $trimmed[$trimmed.Count - 1]
    if ($lastPart.ToLower() -eq "system") {
        $revArr = @()
        for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
            $revArr += $trimmed[$i]
        }
        $newVal = [string]::Join(", ", $revArr)
        $cell.Value = $newVal
    }
}
